$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the closing "fin" row (old row 97), shifting it down to row 100
$ws.Rows("97:99").Insert()

# New "equadiffs" domain rows
$ws.Range("A97").Value = "equadiffs"
$ws.Range("B97").Value = "EQD-004"
$ws.Range("C97").Value = "Equation différentielles d'ordre 1"

$ws.Range("A98").Value = "equadiffs"
$ws.Range("B98").Value = "EQD-005"
$ws.Range("C98").Value = "Equation différentielles d'ordre 2"

$ws.Range("A99").Value = "equadiffs"
$ws.Range("B99").Value = "EQD-006"
$ws.Range("C99").Value = "Cinétique chimique"

# Match the workbook's new selection (was C97, now C100 where "fin" landed)
$null = $ws.Range("C100").Select()
